# Add two new user rows to the "Users" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "vinit"
$ws.Range("C2").Value = "scrypt:32768:8:1`$3YqEzykj0MTLFiSJ`$5320671aedef4422c96ec1a1dcb161a5f997c13be9e432d6a8e39e2e5743698b4b3673d4ea5777ebb26f04c063d19152e252cc55986a070f23e396c173f26958"
$ws.Range("D2").Value = "vinit.j6666@gmail.com"
$ws.Range("E2").Value = "Seller"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "vinit.jadhav"
$ws.Range("C3").Value = "scrypt:32768:8:1`$bcv8YdL86dxDloMF`$e659e065445662d7f7a0c59e46c6b32f3244be4bc9a47a25c73073f1701a7a47efc173ecc0e16206424cce9ad5d53b379cba51aaa6c3492f1bdeecdd2d20c166"
$ws.Range("D3").Value = "vinit.jadhav.vj123@gmail.com"
$ws.Range("E3").Value = "Seller"
